$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'330.73"
$ws.Range("E2").Value = "'0.24%"
$ws.Range("D3").Value = "'43.91"
$ws.Range("E3").Value = "'6.76%"
$ws.Range("D4").Value = "'5.972"
$ws.Range("E4").Value = "'4.83%"
$ws.Range("E5").Value = "'1.58%"
$ws.Range("D6").Value = "'8.788"
$ws.Range("E6").Value = "'0.62%"
$ws.Range("E7").Value = "'-0.82%"
$ws.Range("D8").Value = "'1.973"
$ws.Range("E8").Value = "'-4.86%"
$ws.Range("D9").Value = "'2.915"
$ws.Range("E9").Value = "'-1.58%"
$ws.Range("D10").Value = "'0.9319"
$ws.Range("E10").Value = "'0.64%"
$ws.Range("D11").Value = "'0.1256"
$ws.Range("E11").Value = "'0.10%"
$ws.Range("D12").Value = "'0.1955"
$ws.Range("E12").Value = "'-0.27%"
$ws.Range("D13").Value = "'0.09557"
$ws.Range("E13").Value = "'1.57%"
$ws.Range("E14").Value = "'7.70%"
$ws.Range("D15").Value = "'0.1065"
$ws.Range("E15").Value = "'0.98%"
$ws.Range("D16").Value = "'0.001307"
$ws.Range("E16").Value = "'-0.53%"
$ws.Range("D17").Value = "'0.006018"
$ws.Range("E17").Value = "'-3.62%"
$ws.Range("D18").Value = "'3.533"
$ws.Range("E18").Value = "'3.44%"
$ws.Range("D20").Value = "'9.085"
$ws.Range("E20").Value = "'9.13%"
$ws.Range("E21").Value = "'-0.50%"
$ws.Range("D22").Value = "'0.2603"
$ws.Range("E22").Value = "'-1.92%"
$ws.Range("E23").Value = "'-0.39%"
$ws.Range("D24").Value = "'0.001245"
$ws.Range("E24").Value = "'-1.91%"
$ws.Range("D25").Value = "'0.004397"
$ws.Range("E25").Value = "'1.54%"
$ws.Range("D26").Value = "'0.0001193"
$ws.Range("E26").Value = "'0.90%"
$ws.Range("D27").Value = "'0.0003999"
$ws.Range("E27").Value = "'0.15%"
$ws.Range("D39").Value = "'0.02822"
$ws.Range("E39").Value = "'1.09%"
$ws.Range("D40").Value = "'0.05632"
$ws.Range("E40").Value = "'2.69%"
$ws.Range("D41").Value = "'0.007883"
$ws.Range("E41").Value = "'2.70%"
$ws.Range("E42").Value = "'0.33%"
$ws.Range("D43").Value = "'0.009061"
$ws.Range("E43").Value = "'-3.94%"
$ws.Range("D44").Value = "'0.002115"
$ws.Range("E44").Value = "'-0.89%"
$ws.Range("D45").Value = "'0.008768"
$ws.Range("E45").Value = "'-20.45%"
$ws.Range("D46").Value = "'0.00007315"
$ws.Range("E46").Value = "'6.40%"
$ws.Range("E47").Value = "'0.05%"
$ws.Range("D48").Value = "'0.003597"
$ws.Range("E48").Value = "'11.21%"
$ws.Range("D49").Value = "'0.002283"
$ws.Range("E49").Value = "'-0.02%"
$ws.Range("D50").Value = "'0.00002105"
$ws.Range("E50").Value = "'0.05%"
$ws.Range("D51").Value = "'0.0002005"
$ws.Range("E51").Value = "'0.05%"
